# Update market-price derived columns (H:N) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 175
$ws.Range("I9").Value = 178.44444
$ws.Range("J9").Value = 144
$ws.Range("K9").Value = 178.44444
$ws.Range("L9").Value = 144
$ws.Range("M9").Value = -9.444439999999986
$ws.Range("N9").Value = -482
$ws.Range("H94").Value = 12700
$ws.Range("I94").Value = 2266.6667
$ws.Range("J94").Value = 44000
$ws.Range("K94").Value = 2266.6667
$ws.Range("L94").Value = 44000
$ws.Range("M94").Value = -1815.6667
$ws.Range("N94").Value = -44902
$ws.Range("H125").Value = 3406.125
$ws.Range("I125").Value = 2087.25
$ws.Range("J125").Value = 4725
$ws.Range("K125").Value = 18785.25
$ws.Range("L125").Value = 42525
$ws.Range("M125").Value = -16325.25
$ws.Range("N125").Value = -47445
$ws.Range("H137").Value = 47620770
$ws.Range("I137").Value = 83333890
$ws.Range("J137").Value = 3273.3333
$ws.Range("K137").Value = 250001670
$ws.Range("L137").Value = 9819.999899999999
$ws.Range("M137").Value = -249999120
$ws.Range("N137").Value = -14919.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2443.3914
$ws.Range("I45").Value = 2395.3635
$ws.Range("K45").Value = 2395.3635
$ws.Range("M45").Value = -2018.3635
$ws.Range("H132").Value = 5990.0625
$ws.Range("I132").Value = 4725
$ws.Range("J132").Value = 8098.5
$ws.Range("K132").Value = 14175
$ws.Range("L132").Value = 24295.5
$ws.Range("M132").Value = -11645
$ws.Range("N132").Value = -29355.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1536.3939
$ws.Range("I86").Value = 1434.579
$ws.Range("J86").Value = 1674.5714
$ws.Range("K86").Value = 1434.579
$ws.Range("L86").Value = 1674.5714
$ws.Range("M86").Value = -311.579
$ws.Range("N86").Value = -3920.5714
$ws.Range("H89").Value = 1536.3939
$ws.Range("I89").Value = 1434.579
$ws.Range("J89").Value = 1674.5714
$ws.Range("K89").Value = 7172.895
$ws.Range("L89").Value = 8372.857
$ws.Range("M89").Value = -1556.895
$ws.Range("N89").Value = -19604.857
$ws.Range("H107").Value = 2374.3872
$ws.Range("I107").Value = 2122.4443
$ws.Range("J107").Value = 4075
$ws.Range("K107").Value = 2122.4443
$ws.Range("L107").Value = 4075
$ws.Range("M107").Value = -202.4443000000001
$ws.Range("N107").Value = -7915

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4086.3044
$ws.Range("J31").Value = 5642.5713
$ws.Range("L31").Value = 5642.5713
$ws.Range("N31").Value = -6232.5713
$ws.Range("H34").Value = 4086.3044
$ws.Range("J34").Value = 5642.5713
$ws.Range("L34").Value = 5642.5713
$ws.Range("N34").Value = -6046.5713
$ws.Range("H86").Value = 507490.25
$ws.Range("I86").Value = 9979.5
$ws.Range("K86").Value = 9979.5
$ws.Range("M86").Value = -8856.5
$ws.Range("H89").Value = 507490.25
$ws.Range("I89").Value = 9979.5
$ws.Range("K89").Value = 49897.5
$ws.Range("M89").Value = -44281.5
$ws.Range("H94").Value = 1325.6154
$ws.Range("I94").Value = 973.25
$ws.Range("K94").Value = 973.25
$ws.Range("M94").Value = -522.25
$ws.Range("H99").Value = 2421.8333
$ws.Range("I99").Value = 2032
$ws.Range("K99").Value = 2032
$ws.Range("M99").Value = -534
$ws.Range("H126").Value = 2421.8333
$ws.Range("I126").Value = 2032
$ws.Range("K126").Value = 6096
$ws.Range("M126").Value = -3626

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1273.0605
$ws.Range("I113").Value = 849.75
$ws.Range("J113").Value = 1331.4482
$ws.Range("K113").Value = 2549.25
$ws.Range("L113").Value = 3994.3446
$ws.Range("M113").Value = -379.25
$ws.Range("N113").Value = -8334.3446

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 964.3333
$ws.Range("J97").Value = 1075.375
$ws.Range("L97").Value = 1075.375
$ws.Range("N97").Value = -2067.375
$ws.Range("H102").Value = 4999.3335
$ws.Range("I102").Value = 4999.3335
$ws.Range("K102").Value = 4999.3335
$ws.Range("M102").Value = -3377.3335
$ws.Range("H126").Value = 2255
$ws.Range("I126").Value = 2255
$ws.Range("K126").Value = 6765
$ws.Range("M126").Value = -4295

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4737.696
$ws.Range("I7").Value = 3186.5881
$ws.Range("J7").Value = 9132.5
$ws.Range("K7").Value = 3186.5881
$ws.Range("L7").Value = 9132.5
$ws.Range("M7").Value = -3074.5881
$ws.Range("N7").Value = -9356.5
$ws.Range("H40").Value = 3182.04
$ws.Range("I40").Value = 3318.348
$ws.Range("K40").Value = 3318.348
$ws.Range("M40").Value = -3182.348
$ws.Range("H82").Value = 2115.7
$ws.Range("I82").Value = 1237.8
$ws.Range("J82").Value = 2993.6
$ws.Range("K82").Value = 1237.8
$ws.Range("L82").Value = 2993.6
$ws.Range("M82").Value = -876.8
$ws.Range("N82").Value = -3715.6
$ws.Range("H85").Value = 2115.7
$ws.Range("I85").Value = 1237.8
$ws.Range("J85").Value = 2993.6
$ws.Range("K85").Value = 1237.8
$ws.Range("L85").Value = 2993.6
$ws.Range("M85").Value = 10.20000000000005
$ws.Range("N85").Value = -5489.6
$ws.Range("H122").Value = 4895.074
$ws.Range("I122").Value = 4458.95
$ws.Range("J122").Value = 6141.143
$ws.Range("K122").Value = 13376.85
$ws.Range("L122").Value = 18423.429
$ws.Range("M122").Value = -10926.85
$ws.Range("N122").Value = -23323.429
$ws.Range("H126").Value = 4737.696
$ws.Range("I126").Value = 3186.5881
$ws.Range("J126").Value = 9132.5
$ws.Range("K126").Value = 9559.764299999999
$ws.Range("L126").Value = 27397.5
$ws.Range("M126").Value = -7089.764299999999
$ws.Range("N126").Value = -32337.5
$ws.Range("H136").Value = 3480.2703
$ws.Range("I136").Value = 1861.2858
$ws.Range("K136").Value = 5583.857400000001
$ws.Range("M136").Value = -3033.857400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 91144690
$ws.Range("I4").Value = 630748.75
$ws.Range("J4").Value = 142866940
$ws.Range("K4").Value = 630748.75
$ws.Range("L4").Value = 142866940
$ws.Range("M4").Value = -630635.75
$ws.Range("N4").Value = -142867166
$ws.Range("H122").Value = 3789.6667
$ws.Range("I122").Value = 3550.2222
$ws.Range("J122").Value = 6663
$ws.Range("K122").Value = 10650.6666
$ws.Range("L122").Value = 19989
$ws.Range("M122").Value = -8200.6666
$ws.Range("N122").Value = -24889
